# Daily attendance processing - 2026-01-14 12:54:22
# Normalize the "Recorded By" (column G) entries so that "System" is
# listed before the recorder's email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($i = 1; $i -le $rowCount; $i++) {
    $cell = $ws.Cells.Item($i, 7)   # Column G
    if ($cell.Text -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}
